# Weekly crime-data refresh for 19th Precinct CompStat report
# (Volume 30, Number 5 - week covering 1/30/2023 through 2/5/2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header text (week number and date range) ---
$ws.Range("A8").Value2 = "Volume 30   Number  5"
$ws.Range("C9").Value2 = "Report Covering the Week  1/30/2023  Through  2/5/2023"

# --- Update CompStat data table (rows 15-30) ---
# Note: some cells switch between a numeric value and the text placeholders
# used elsewhere in the sheet for "not applicable" (shown as "0") or an
# undefined percentage change (shown as "***.*"). Style + shared-text for
# those placeholder cells are copied from untouched reference cells that
# already use the same style (row 14, and F15 for the plain numeric style),
# then overwritten with the correct value where the target cell is numeric.
# Row 15
$ws.Range("C14").Copy($ws.Range("C15"))
# Row 16
$ws.Range("C16").Value2 = 2
$ws.Range("F15").Copy($ws.Range("D16"))
$ws.Range("D16").Value2 = 5
$ws.Range("M14").Copy($ws.Range("E16"))
$ws.Range("E16").Value2 = -60
$ws.Range("F16").Value2 = 13
$ws.Range("H16").Value2 = 0
$ws.Range("I16").Value2 = 17
$ws.Range("J16").Value2 = 19
$ws.Range("K16").Value2 = -10.526315789473
$ws.Range("L16").Value2 = 21.428571428571
$ws.Range("M16").Value2 = 21.428571428571
$ws.Range("N16").Value2 = -87.681159420289
# Row 17
$ws.Range("C17").Value2 = 5
$ws.Range("D17").Value2 = 2
$ws.Range("E17").Value2 = 150
$ws.Range("F17").Value2 = 17
$ws.Range("G17").Value2 = 15
$ws.Range("H17").Value2 = 13.333333333333
$ws.Range("I17").Value2 = 21
$ws.Range("J17").Value2 = 16
$ws.Range("K17").Value2 = 31.25
$ws.Range("L17").Value2 = 90.909090909090
$ws.Range("M17").Value2 = 250
$ws.Range("N17").Value2 = 16.666666666666
# Row 18
$ws.Range("C18").Value2 = 2
$ws.Range("D18").Value2 = 4
$ws.Range("E18").Value2 = -50
$ws.Range("F18").Value2 = 13
$ws.Range("H18").Value2 = -53.571428571428
$ws.Range("I18").Value2 = 21
$ws.Range("J18").Value2 = 35
$ws.Range("K18").Value2 = -40
$ws.Range("L18").Value2 = -41.666666666666
$ws.Range("M18").Value2 = -48.780487804878
$ws.Range("N18").Value2 = -92.250922509225
# Row 19
$ws.Range("C19").Value2 = 25
$ws.Range("E19").Value2 = -24.242424242424
$ws.Range("F19").Value2 = 114
$ws.Range("G19").Value2 = 108
$ws.Range("H19").Value2 = 5.555555555555
$ws.Range("I19").Value2 = 136
$ws.Range("J19").Value2 = 132
$ws.Range("K19").Value2 = 3.030303030303
$ws.Range("L19").Value2 = 34.653465346534
$ws.Range("M19").Value2 = 12.396694214876
$ws.Range("N19").Value2 = -63.243243243243
# Row 20
$ws.Range("C20").Value2 = 2
$ws.Range("D20").Value2 = 2
$ws.Range("E20").Value2 = 0
$ws.Range("F20").Value2 = 7
$ws.Range("G20").Value2 = 8
$ws.Range("H20").Value2 = -12.5
$ws.Range("I20").Value2 = 11
$ws.Range("J20").Value2 = 12
$ws.Range("K20").Value2 = -8.333333333333
$ws.Range("L20").Value2 = 10
$ws.Range("M20").Value2 = 175
$ws.Range("N20").Value2 = -96.978021978022
# Row 21
$ws.Range("C21").Value2 = 36
$ws.Range("D21").Value2 = 46
$ws.Range("E21").Value2 = -21.739130434782
$ws.Range("F21").Value2 = 166
$ws.Range("G21").Value2 = 173
$ws.Range("H21").Value2 = -4.046242774566
$ws.Range("I21").Value2 = 208
$ws.Range("J21").Value2 = 215
$ws.Range("K21").Value2 = -3.255813953488
$ws.Range("L21").Value2 = 20.930232558139
$ws.Range("M21").Value2 = 10.638297872340
$ws.Range("N21").Value2 = -82.145922746781
# Row 22
$ws.Range("F15").Copy($ws.Range("D22"))
$ws.Range("D22").Value2 = 3
$ws.Range("M14").Copy($ws.Range("E22"))
$ws.Range("E22").Value2 = -100
$ws.Range("G22").Value2 = 5
$ws.Range("J22").Value2 = 6
# Row 23
$ws.Range("C14").Copy($ws.Range("C23"))
$ws.Range("C14").Copy($ws.Range("D23"))
$ws.Range("E14").Copy($ws.Range("E23"))
$ws.Range("L23").Value2 = -66.666666666666
# Row 24
$ws.Range("C24").Value2 = 63
$ws.Range("D24").Value2 = 48
$ws.Range("E24").Value2 = 31.25
$ws.Range("F24").Value2 = 261
$ws.Range("G24").Value2 = 220
$ws.Range("H24").Value2 = 18.636363636363
$ws.Range("I24").Value2 = 314
$ws.Range("J24").Value2 = 272
$ws.Range("K24").Value2 = 15.441176470588
$ws.Range("L24").Value2 = 35.344827586206
$ws.Range("M24").Value2 = 91.463414634146
# Row 25
$ws.Range("C25").Value2 = 8
$ws.Range("E25").Value2 = 33.333333333333
$ws.Range("F25").Value2 = 31
$ws.Range("G25").Value2 = 29
$ws.Range("H25").Value2 = 6.896551724137
$ws.Range("I25").Value2 = 37
$ws.Range("J25").Value2 = 40
$ws.Range("K25").Value2 = -7.5
$ws.Range("L25").Value2 = 68.181818181818
$ws.Range("M25").Value2 = 5.714285714285
# Row 26
$ws.Range("C26").Value2 = 1
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("F26").Value2 = 4
$ws.Range("H26").Value2 = 100
$ws.Range("I26").Value2 = 4
$ws.Range("K26").Value2 = 100
# Row 27
$ws.Range("C27").Value2 = 1
$ws.Range("F15").Copy($ws.Range("D27"))
$ws.Range("D27").Value2 = 3
$ws.Range("M14").Copy($ws.Range("E27"))
$ws.Range("E27").Value2 = -66.666666666666
$ws.Range("F27").Value2 = 9
$ws.Range("G27").Value2 = 5
$ws.Range("H27").Value2 = 80
$ws.Range("I27").Value2 = 12
$ws.Range("J27").Value2 = 12
$ws.Range("K27").Value2 = 0
$ws.Range("L27").Value2 = 140
# Row 28
$ws.Range("C14").Copy($ws.Range("G28"))
$ws.Range("E14").Copy($ws.Range("H28"))
# Row 29
$ws.Range("C14").Copy($ws.Range("G29"))
$ws.Range("E14").Copy($ws.Range("H29"))
# Row 30
$ws.Range("C14").Copy($ws.Range("F30"))
$ws.Range("H30").Value2 = -100
$ws.Range("M14").Copy($ws.Range("L30"))
$ws.Range("L30").Value2 = 0
